$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 325 (shifts existing row 325..348 down to 326..349)
$ws.Rows(325).Insert()

# Populate the newly inserted row 325 with the new record.
# Unchanged columns mirror the record that used to occupy row 325 (now row 326).
$ws.Range("A325").Value = 11
$ws.Range("B325").Value = "Vega Monumental Concepción"
$ws.Range("C325").Value = "Bíobío"
$ws.Range("D325").Value = "2023-01-05"
$ws.Range("E325").Value = 8
$ws.Range("F325").Value = 100112009
$ws.Range("G325").Value = "Acelga"
$ws.Range("H325").Value = "Sin especificar"
$ws.Range("I325").Value = "Primera"
$ws.Range("J325").Value = 350
$ws.Range("K325").Value = 600
$ws.Range("L325").Value = 650
$ws.Range("M325").Value = 629
$ws.Range("N325").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O325").Value = "Región de Ñuble"
$ws.Range("P325").Value = 629
$ws.Range("Q325").Value = 1
$ws.Range("R325").Value = "Hortaliza"
